$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.75"
$ws.Range("D2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.92"
$ws.Range("D3").NumberFormat = "General"

$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E4").Value = "3LEOLEO"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3.499"
$ws.Range("D4").NumberFormat = "General"

$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "5.171"
$ws.Range("D5").NumberFormat = "General"

$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.05725"
$ws.Range("D6").NumberFormat = "General"

$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.486"
$ws.Range("D7").NumberFormat = "General"

$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("E8").Value = "7GateTokenGT"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.063"
$ws.Range("D8").NumberFormat = "General"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8103"
$ws.Range("D9").NumberFormat = "General"

$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("E10").Value = "9FTXTokenFTT"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8380"
$ws.Range("D10").NumberFormat = "General"

$ws.Range("B11").Value = "One"
$ws.Range("C11").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E11").Value = "10OneONE"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0006000"
$ws.Range("D11").NumberFormat = "General"

$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E12").Value = "11WazirXWRX"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1336"
$ws.Range("D12").NumberFormat = "General"

$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.06967"
$ws.Range("D13").NumberFormat = "General"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02820"
$ws.Range("D14").NumberFormat = "General"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09378"
$ws.Range("D15").NumberFormat = "General"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001510"
$ws.Range("D16").NumberFormat = "General"

$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006241"
$ws.Range("D17").NumberFormat = "General"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3193"
$ws.Range("D19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03235"
$ws.Range("D20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1300"
$ws.Range("D21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.741"
$ws.Range("D22").NumberFormat = "General"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04669"
$ws.Range("D23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1328"
$ws.Range("D24").NumberFormat = "General"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004249"
$ws.Range("D26").NumberFormat = "General"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009700"
$ws.Range("D27").NumberFormat = "General"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03627"
$ws.Range("D40").NumberFormat = "General"

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006320"
$ws.Range("D41").NumberFormat = "General"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1048"
$ws.Range("D42").NumberFormat = "General"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003000"
$ws.Range("D43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007335"
$ws.Range("D44").NumberFormat = "General"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005279"
$ws.Range("D45").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1900"
$ws.Range("D47").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D49").NumberFormat = "General"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("D50").NumberFormat = "General"
